$d = $word.ActiveDocument

$null = $d.Content.Find.Execute("69+29=98", $true, $false, $false, $false, $false, $true, 1, $false, "60-38=22", 2)
$null = $d.Content.Find.Execute("58+26=84", $true, $false, $false, $false, $false, $true, 1, $false, "49+19=68", 2)
$null = $d.Content.Find.Execute("88-39=49", $true, $false, $false, $false, $false, $true, 1, $false, "72+9=81", 2)
$null = $d.Content.Find.Execute("71-64=7", $true, $false, $false, $false, $false, $true, 1, $false, "91-58=33", 2)
$null = $d.Content.Find.Execute("73-37=36", $true, $false, $false, $false, $false, $true, 1, $false, "61-47=14", 2)
$null = $d.Content.Find.Execute("84-46=38", $true, $false, $false, $false, $false, $true, 1, $false, "84-9=75", 2)
$null = $d.Content.Find.Execute("26+38=64", $true, $false, $false, $false, $false, $true, 1, $false, "61-27=34", 2)
$null = $d.Content.Find.Execute("68+19=87", $true, $false, $false, $false, $false, $true, 1, $false, "93-4=89", 2)
$null = $d.Content.Find.Execute("49+37=86", $true, $false, $false, $false, $false, $true, 1, $false, "56+9=65", 2)
$null = $d.Content.Find.Execute("27+56=83", $true, $false, $false, $false, $false, $true, 1, $false, "72+19=91", 2)
$null = $d.Content.Find.Execute("83-4=79", $true, $false, $false, $false, $false, $true, 1, $false, "61-6=55", 2)
$null = $d.Content.Find.Execute("40-39=1", $true, $false, $false, $false, $false, $true, 1, $false, "7+74=81", 2)
$null = $d.Content.Find.Execute("53-18=35", $true, $false, $false, $false, $false, $true, 1, $false, "97-39=58", 2)
$null = $d.Content.Find.Execute("70-7=63", $true, $false, $false, $false, $false, $true, 1, $false, "35-28=7", 2)
$null = $d.Content.Find.Execute("19+37=56", $true, $false, $false, $false, $false, $true, 1, $false, "15+77=92", 2)
$null = $d.Content.Find.Execute("81-3=78", $true, $false, $false, $false, $false, $true, 1, $false, "37-9=28", 2)
$null = $d.Content.Find.Execute("8+88=96", $true, $false, $false, $false, $false, $true, 1, $false, "8+87=95", 2)
$null = $d.Content.Find.Execute("19+42=61", $true, $false, $false, $false, $false, $true, 1, $false, "57+39=96", 2)
$null = $d.Content.Find.Execute("15+67=82", $true, $false, $false, $false, $false, $true, 1, $false, "95-89=6", 2)
$null = $d.Content.Find.Execute("65+19=84", $true, $false, $false, $false, $false, $true, 1, $false, "84-76=8", 2)
$null = $d.Content.Find.Execute("57+34=91", $true, $false, $false, $false, $false, $true, 1, $false, "63-58=5", 2)
$null = $d.Content.Find.Execute("6+75=81", $true, $false, $false, $false, $false, $true, 1, $false, "19+34=53", 2)
$null = $d.Content.Find.Execute("82-7=75", $true, $false, $false, $false, $false, $true, 1, $false, "7+6=13", 2)
$null = $d.Content.Find.Execute("24+7=31", $true, $false, $false, $false, $false, $true, 1, $false, "4+79=83", 2)
$null = $d.Content.Find.Execute("32-5=27", $true, $false, $false, $false, $false, $true, 1, $false, "67+18=85", 2)
$null = $d.Content.Find.Execute("2+59=61", $true, $false, $false, $false, $false, $true, 1, $false, "77+7=84", 2)
$null = $d.Content.Find.Execute("60-32=28", $true, $false, $false, $false, $false, $true, 1, $false, "90-2=88", 2)
$null = $d.Content.Find.Execute("55+38=93", $true, $false, $false, $false, $false, $true, 1, $false, "34-18=16", 2)
$null = $d.Content.Find.Execute("58+7=65", $true, $false, $false, $false, $false, $true, 1, $false, "57+15=72", 2)
$null = $d.Content.Find.Execute("64-55=9", $true, $false, $false, $false, $false, $true, 1, $false, "85-67=18", 2)
$null = $d.Content.Find.Execute("87-79=8", $true, $false, $false, $false, $false, $true, 1, $false, "13+49=62", 2)
$null = $d.Content.Find.Execute("74-17=57", $true, $false, $false, $false, $false, $true, 1, $false, "55-36=19", 2)
$null = $d.Content.Find.Execute("19+4=23", $true, $false, $false, $false, $false, $true, 1, $false, "61-43=18", 2)
$null = $d.Content.Find.Execute("72-63=9", $true, $false, $false, $false, $false, $true, 1, $false, "61-28=33", 2)
$null = $d.Content.Find.Execute("46+39=85", $true, $false, $false, $false, $false, $true, 1, $false, "58+36=94", 2)
$null = $d.Content.Find.Execute("53+28=81", $true, $false, $false, $false, $false, $true, 1, $false, "97-38=59", 2)
$null = $d.Content.Find.Execute("42-7=35", $true, $false, $false, $false, $false, $true, 1, $false, "73-19=54", 2)
$null = $d.Content.Find.Execute("93-15=78", $true, $false, $false, $false, $false, $true, 1, $false, "8+39=47", 2)
$null = $d.Content.Find.Execute("78+3=81", $true, $false, $false, $false, $false, $true, 1, $false, "15+38=53", 2)
$null = $d.Content.Find.Execute("18+63=81", $true, $false, $false, $false, $false, $true, 1, $false, "47+9=56", 2)
$null = $d.Content.Find.Execute("42+19=61", $true, $false, $false, $false, $false, $true, 1, $false, "62-57=5", 2)
$null = $d.Content.Find.Execute("91-82=9", $true, $false, $false, $false, $false, $true, 1, $false, "77-38=39", 2)
$null = $d.Content.Find.Execute("19+76=95", $true, $false, $false, $false, $false, $true, 1, $false, "32+19=51", 2)
$null = $d.Content.Find.Execute("58-9=49", $true, $false, $false, $false, $false, $true, 1, $false, "68+28=96", 2)
$null = $d.Content.Find.Execute("8+55=63", $true, $false, $false, $false, $false, $true, 1, $false, "70-59=11", 2)
$null = $d.Content.Find.Execute("85-77=8", $true, $false, $false, $false, $false, $true, 1, $false, "43-4=39", 2)
$null = $d.Content.Find.Execute("25+38=63", $true, $false, $false, $false, $false, $true, 1, $false, "66+17=83", 2)
$null = $d.Content.Find.Execute("80-64=16", $true, $false, $false, $false, $false, $true, 1, $false, "69+13=82", 2)
$null = $d.Content.Find.Execute("48+8=56", $true, $false, $false, $false, $false, $true, 1, $false, "56+29=85", 2)
$null = $d.Content.Find.Execute("53-16=37", $true, $false, $false, $false, $false, $true, 1, $false, "98-29=69", 2)
$null = $d.Content.Find.Execute("73-7=66", $true, $false, $false, $false, $false, $true, 1, $false, "78+13=91", 2)
$null = $d.Content.Find.Execute("93-45=48", $true, $false, $false, $false, $false, $true, 1, $false, "90-52=38", 2)
$null = $d.Content.Find.Execute("44-6=38", $true, $false, $false, $false, $false, $true, 1, $false, "19+29=48", 2)
$null = $d.Content.Find.Execute("82-35=47", $true, $false, $false, $false, $false, $true, 1, $false, "18+66=84", 2)
$null = $d.Content.Find.Execute("70-64=6", $true, $false, $false, $false, $false, $true, 1, $false, "35+59=94", 2)
$null = $d.Content.Find.Execute("44+28=72", $true, $false, $false, $false, $false, $true, 1, $false, "43+19=62", 2)
$null = $d.Content.Find.Execute("7+49=56", $true, $false, $false, $false, $false, $true, 1, $false, "59+26=85", 2)
$null = $d.Content.Find.Execute("50-32=18", $true, $false, $false, $false, $false, $true, 1, $false, "61-13=48", 2)
$null = $d.Content.Find.Execute("46+5=51", $true, $false, $false, $false, $false, $true, 1, $false, "17+66=83", 2)
$null = $d.Content.Find.Execute("46+8=54", $true, $false, $false, $false, $false, $true, 1, $false, "51-15=36", 2)
$null = $d.Content.Find.Execute("27+67=94", $true, $false, $false, $false, $false, $true, 1, $false, "59+6=65", 2)
$null = $d.Content.Find.Execute("65-17=48", $true, $false, $false, $false, $false, $true, 1, $false, "4+27=31", 2)
$null = $d.Content.Find.Execute("48+15=63", $true, $false, $false, $false, $false, $true, 1, $false, "24+57=81", 2)
$null = $d.Content.Find.Execute("28+7=35", $true, $false, $false, $false, $false, $true, 1, $false, "36+57=93", 2)
$null = $d.Content.Find.Execute("75-47=28", $true, $false, $false, $false, $false, $true, 1, $false, "69+6=75", 2)
$null = $d.Content.Find.Execute("66-19=47", $true, $false, $false, $false, $false, $true, 1, $false, "72-54=18", 2)
$null = $d.Content.Find.Execute("29+28=57", $true, $false, $false, $false, $false, $true, 1, $false, "58+8=66", 2)
$null = $d.Content.Find.Execute("19+77=96", $true, $false, $false, $false, $false, $true, 1, $false, "36+18=54", 2)
$null = $d.Content.Find.Execute("33-17=16", $true, $false, $false, $false, $false, $true, 1, $false, "56+27=83", 2)
$null = $d.Content.Find.Execute("9+25=34", $true, $false, $false, $false, $false, $true, 1, $false, "84-76=8", 2)
$null = $d.Content.Find.Execute("75-6=69", $true, $false, $false, $false, $false, $true, 1, $false, "28+69=97", 2)
$null = $d.Content.Find.Execute("18+44=62", $true, $false, $false, $false, $false, $true, 1, $false, "35-6=29", 2)
$null = $d.Content.Find.Execute("38+59=97", $true, $false, $false, $false, $false, $true, 1, $false, "23+18=41", 2)
$null = $d.Content.Find.Execute("22+59=81", $true, $false, $false, $false, $false, $true, 1, $false, "38+58=96", 2)
$null = $d.Content.Find.Execute("19+46=65", $true, $false, $false, $false, $false, $true, 1, $false, "19+55=74", 2)
$null = $d.Content.Find.Execute("38+19=57", $true, $false, $false, $false, $false, $true, 1, $false, "27+27=54", 2)
$null = $d.Content.Find.Execute("96-67=29", $true, $false, $false, $false, $false, $true, 1, $false, "15+36=51", 2)
$null = $d.Content.Find.Execute("29+63=92", $true, $false, $false, $false, $false, $true, 1, $false, "51-7=44", 2)
$null = $d.Content.Find.Execute("26+28=54", $true, $false, $false, $false, $false, $true, 1, $false, "28+63=91", 2)
$null = $d.Content.Find.Execute("92-5=87", $true, $false, $false, $false, $false, $true, 1, $false, "61-26=35", 2)
$null = $d.Content.Find.Execute("49+49=98", $true, $false, $false, $false, $false, $true, 1, $false, "34+28=62", 2)
$null = $d.Content.Find.Execute("73-14=59", $true, $false, $false, $false, $false, $true, 1, $false, "61-37=24", 2)
$null = $d.Content.Find.Execute("75-8=67", $true, $false, $false, $false, $false, $true, 1, $false, "90-5=85", 2)
$null = $d.Content.Find.Execute("44-18=26", $true, $false, $false, $false, $false, $true, 1, $false, "35+27=62", 2)
$null = $d.Content.Find.Execute("87-49=38", $true, $false, $false, $false, $false, $true, 1, $false, "70-52=18", 2)
$null = $d.Content.Find.Execute("74-48=26", $true, $false, $false, $false, $false, $true, 1, $false, "18+68=86", 2)
$null = $d.Content.Find.Execute("74-49=25", $true, $false, $false, $false, $false, $true, 1, $false, "70-42=28", 2)
$null = $d.Content.Find.Execute("50-1=49", $true, $false, $false, $false, $false, $true, 1, $false, "92-27=65", 2)
$null = $d.Content.Find.Execute("91-46=45", $true, $false, $false, $false, $false, $true, 1, $false, "29+37=66", 2)
$null = $d.Content.Find.Execute("40-38=2", $true, $false, $false, $false, $false, $true, 1, $false, "90-23=67", 2)
$null = $d.Content.Find.Execute("22-17=5", $true, $false, $false, $false, $false, $true, 1, $false, "61-46=15", 2)
$null = $d.Content.Find.Execute("16+65=81", $true, $false, $false, $false, $false, $true, 1, $false, "48+4=52", 2)
$null = $d.Content.Find.Execute("74+19=93", $true, $false, $false, $false, $false, $true, 1, $false, "28+19=47", 2)
$null = $d.Content.Find.Execute("68+24=92", $true, $false, $false, $false, $false, $true, 1, $false, "38+37=75", 2)
$null = $d.Content.Find.Execute("86-9=77", $true, $false, $false, $false, $false, $true, 1, $false, "95-37=58", 2)
$null = $d.Content.Find.Execute("40-8=32", $true, $false, $false, $false, $false, $true, 1, $false, "24+28=52", 2)
$null = $d.Content.Find.Execute("35+37=72", $true, $false, $false, $false, $false, $true, 1, $false, "5+88=93", 2)
$null = $d.Content.Find.Execute("28+18=46", $true, $false, $false, $false, $false, $true, 1, $false, "36+16=52", 2)
$null = $d.Content.Find.Execute("41-25=16", $true, $false, $false, $false, $false, $true, 1, $false, "15+7=22", 2)
$null = $d.Content.Find.Execute("24+47=71", $true, $false, $false, $false, $false, $true, 1, $false, "94-18=76", 2)
